# "updated main GSC export data"
#
# The GSC export adds one more day ("2025-12-17") of data to the bottom of
# the "Chart" sheet's Date/Invalid/Valid table (previously ending at row 73
# with 2025-12-16 / 0 / 32). The new row keeps 0 invalid items and closes
# out with 31 valid items.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Find the first empty row right after the existing data (row 73 -> 74).
$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

# The Date column stores plain text that looks like "yyyy-MM-dd" (not real
# dates). Prefix with an apostrophe so Excel keeps it as text instead of
# auto-converting it to a date serial, matching every other cell in column A.
$ws.Cells.Item($newRow, 1).Value = "'2025-12-17"
$ws.Cells.Item($newRow, 2).Value = 0
$ws.Cells.Item($newRow, 3).Value = 31
